$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 / Row 6: swap match data (MC Alger vs Magra matches on the same date) ---
$ws.Range("F5").Value = 'Magra'
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 'Kabylie'
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 2.84
$ws.Range("K5").Value = '15/09/2023 13:42'
$ws.Range("L5").Value = 3.8
$ws.Range("M5").Value = '16/09/2023 16:12'
$ws.Range("N5").Value = 2.63
$ws.Range("O5").Value = '15/09/2023 13:42'
$ws.Range("P5").Value = 2.84
$ws.Range("Q5").Value = '16/09/2023 15:03'
$ws.Range("R5").Value = 2.72
$ws.Range("S5").Value = '15/09/2023 13:42'
$ws.Range("T5").Value = 2.26
$ws.Range("U5").Value = '16/09/2023 16:12'
$ws.Range("V5").Value = 'https://www.betexplorer.com/football/algeria/ligue-1/magra-kabylie/YFXa8c8H/'
$ws.Range("F6").Value = 'MC Alger'
$ws.Range("G6").Value = 4
$ws.Range("H6").Value = 'Ben Aknoun'
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 1.4
$ws.Range("K6").Value = '16/09/2023 03:43'
$ws.Range("L6").Value = 1.3
$ws.Range("M6").Value = '16/09/2023 10:40'
$ws.Range("N6").Value = 4.19
$ws.Range("O6").Value = '16/09/2023 03:43'
$ws.Range("P6").Value = 4.81
$ws.Range("Q6").Value = '16/09/2023 16:47'
$ws.Range("R6").Value = 8.529999999999999
$ws.Range("S6").Value = '16/09/2023 03:43'
$ws.Range("T6").Value = 12.64
$ws.Range("U6").Value = '16/09/2023 16:47'
$ws.Range("V6").Value = 'https://www.betexplorer.com/football/algeria/ligue-1/mc-alger-es-ben-aknoun/WjyqCu9h/'

# --- Rows 16-18: rotate match data (US Souf / Magra / Khenchela matches on the same date) ---
$ws.Range("F16").Value = 'Magra'
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 'Biskra'
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 1.98
$ws.Range("K16").Value = '28/09/2023 04:12'
$ws.Range("L16").Value = 2.1
$ws.Range("M16").Value = '29/09/2023 16:44'
$ws.Range("N16").Value = 2.89
$ws.Range("O16").Value = '28/09/2023 04:12'
$ws.Range("P16").Value = 2.77
$ws.Range("Q16").Value = '29/09/2023 16:44'
$ws.Range("R16").Value = 3.97
$ws.Range("S16").Value = '28/09/2023 04:12'
$ws.Range("T16").Value = 4.54
$ws.Range("U16").Value = '29/09/2023 16:44'
$ws.Range("V16").Value = 'https://www.betexplorer.com/football/algeria/ligue-1/magra-biskra/OKYxGuDl/'
$ws.Range("F17").Value = 'Khenchela'
$ws.Range("G17").Value = 2
$ws.Range("H17").Value = 'Kabylie'
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = 2.63
$ws.Range("K17").Value = '28/09/2023 04:12'
$ws.Range("L17").Value = 2.05
$ws.Range("M17").Value = '29/09/2023 16:41'
$ws.Range("N17").Value = 2.62
$ws.Range("O17").Value = '28/09/2023 04:12'
$ws.Range("P17").Value = 2.75
$ws.Range("Q17").Value = '29/09/2023 16:41'
$ws.Range("R17").Value = 3.02
$ws.Range("S17").Value = '28/09/2023 04:12'
$ws.Range("T17").Value = 4.88
$ws.Range("U17").Value = '29/09/2023 16:27'
$ws.Range("V17").Value = 'https://www.betexplorer.com/football/algeria/ligue-1/khenchela-kabylie/pUZYGLcr/'
$ws.Range("F18").Value = 'US Souf'
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 'Oran'
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 2.49
$ws.Range("K18").Value = '28/09/2023 19:27'
$ws.Range("L18").Value = 2.14
$ws.Range("M18").Value = '29/09/2023 13:29'
$ws.Range("N18").Value = 2.88
$ws.Range("O18").Value = '28/09/2023 19:27'
$ws.Range("P18").Value = 2.74
$ws.Range("Q18").Value = '29/09/2023 14:49'
$ws.Range("R18").Value = 3.18
$ws.Range("S18").Value = '28/09/2023 19:27'
$ws.Range("T18").Value = 4.43
$ws.Range("U18").Value = '29/09/2023 15:47'
$ws.Range("V18").Value = 'https://www.betexplorer.com/football/algeria/ligue-1/us-souf-oran/6qOsFaSf/'

# --- Row 31 / Row 32: swap match data (Paradou vs US Souf matches on the same date) ---
$ws.Range("F31").Value = 'US Souf'
$ws.Range("G31").Value = 3
$ws.Range("H31").Value = 'Constantine'
$ws.Range("I31").Value = 4
$ws.Range("J31").Value = 3.14
$ws.Range("K31").Value = '10/11/2023 06:42'
$ws.Range("L31").Value = 3.21
$ws.Range("M31").Value = '10/11/2023 15:16'
$ws.Range("N31").Value = 2.82
$ws.Range("O31").Value = '10/11/2023 06:42'
$ws.Range("P31").Value = 2.86
$ws.Range("Q31").Value = '10/11/2023 13:35'
$ws.Range("R31").Value = 2.5
$ws.Range("S31").Value = '10/11/2023 06:42'
$ws.Range("T31").Value = 2.53
$ws.Range("U31").Value = '10/11/2023 15:16'
$ws.Range("V31").Value = 'https://www.betexplorer.com/football/algeria/ligue-1/us-souf-constantine/6mEJaZvD/'
$ws.Range("F32").Value = 'Paradou'
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = 'Oran'
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1.48
$ws.Range("K32").Value = '19/10/2023 04:42'
$ws.Range("L32").Value = 1.46
$ws.Range("M32").Value = '10/11/2023 14:57'
$ws.Range("N32").Value = 3.71
$ws.Range("O32").Value = '19/10/2023 04:42'
$ws.Range("P32").Value = 3.91
$ws.Range("Q32").Value = '10/11/2023 14:57'
$ws.Range("R32").Value = 6.42
$ws.Range("S32").Value = '19/10/2023 04:42'
$ws.Range("T32").Value = 8.550000000000001
$ws.Range("U32").Value = '10/11/2023 14:57'
$ws.Range("V32").Value = 'https://www.betexplorer.com/football/algeria/ligue-1/paradou-oran/ALKA1eA0/'

# --- New rows 58-60: append 3 new matches, copying row 57 formatting first ---
$ws.Range("A57:V57").Copy()
$ws.Range("A58:V60").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A58").Value = 57
$ws.Range("B58").Value = 'algeria'
$ws.Range("C58").Value = 'ligue-1'
$ws.Range("D58").Value = '2023-2024'
$ws.Range("E58").Value = 45262.69791666666
$ws.Range("F58").Value = 'Saoura'
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 'ES Setif'
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 1.88
$ws.Range("K58").Value = '01/12/2023 05:12'
$ws.Range("L58").Value = 1.75
$ws.Range("M58").Value = '02/12/2023 16:40'
$ws.Range("N58").Value = 3.01
$ws.Range("O58").Value = '01/12/2023 05:12'
$ws.Range("P58").Value = 3.32
$ws.Range("Q58").Value = '02/12/2023 16:40'
$ws.Range("R58").Value = 4.3
$ws.Range("S58").Value = '01/12/2023 05:12'
$ws.Range("T58").Value = 5.44
$ws.Range("U58").Value = '02/12/2023 16:40'
$ws.Range("V58").Value = 'https://www.betexplorer.com/football/algeria/ligue-1/saoura-es-setif/Kb2wh3b1/'

$ws.Range("A59").Value = 58
$ws.Range("B59").Value = 'algeria'
$ws.Range("C59").Value = 'ligue-1'
$ws.Range("D59").Value = '2023-2024'
$ws.Range("E59").Value = 45262.69791666666
$ws.Range("F59").Value = 'ASO Chlef'
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = 'US Souf'
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 1.35
$ws.Range("K59").Value = '01/12/2023 05:12'
$ws.Range("L59").Value = 1.32
$ws.Range("M59").Value = '02/12/2023 16:27'
$ws.Range("N59").Value = 4.37
$ws.Range("O59").Value = '01/12/2023 05:12'
$ws.Range("P59").Value = 4.85
$ws.Range("Q59").Value = '02/12/2023 16:27'
$ws.Range("R59").Value = 8.279999999999999
$ws.Range("S59").Value = '01/12/2023 05:12'
$ws.Range("T59").Value = 11.33
$ws.Range("U59").Value = '02/12/2023 16:27'
$ws.Range("V59").Value = 'https://www.betexplorer.com/football/algeria/ligue-1/aso-chlef-us-souf/xzDYhqqe/'

$ws.Range("A60").Value = 59
$ws.Range("B60").Value = 'algeria'
$ws.Range("C60").Value = 'ligue-1'
$ws.Range("D60").Value = '2023-2024'
$ws.Range("E60").Value = 45262.75
$ws.Range("F60").Value = 'Oran'
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 'Kabylie'
$ws.Range("I60").Value = 3
$ws.Range("J60").Value = 2.38
$ws.Range("K60").Value = '01/12/2023 06:13'
$ws.Range("L60").Value = 2.77
$ws.Range("M60").Value = '02/12/2023 17:55'
$ws.Range("N60").Value = 2.48
$ws.Range("O60").Value = '01/12/2023 06:13'
$ws.Range("P60").Value = 2.65
$ws.Range("Q60").Value = '02/12/2023 17:55'
$ws.Range("R60").Value = 3.75
$ws.Range("S60").Value = '01/12/2023 06:13'
$ws.Range("T60").Value = 3.14
$ws.Range("U60").Value = '02/12/2023 17:55'
$ws.Range("V60").Value = 'https://www.betexplorer.com/football/algeria/ligue-1/oran-kabylie/pxfPF6j8/'

